$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2999.6667
$ws.Range("I6").Value = 2999.6667
$ws.Range("K6").Value = 8999.000100000001
$ws.Range("M6").Value = -8887.000100000001
$ws.Range("H9").Value = 3426104.2
$ws.Range("I9").Value = 733
$ws.Range("J9").Value = 7993265.5
$ws.Range("K9").Value = 733
$ws.Range("L9").Value = 7993265.5
$ws.Range("M9").Value = -564
$ws.Range("N9").Value = -7993603.5
$ws.Range("H12").Value = 17843.666
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 17843.666
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 17843.666
$ws.Range("N12").Value = -18183.666
$ws.Range("H17").Value = 6295.7144
$ws.Range("J17").Value = 6295.7144
$ws.Range("L17").Value = 18887.1432
$ws.Range("N17").Value = -19223.1432
$ws.Range("H26").Value = 11166.667
$ws.Range("J26").Value = 12250
$ws.Range("L26").Value = 12250
$ws.Range("N26").Value = -12938
$ws.Range("H28").Value = 213.38461
$ws.Range("I28").Value = 242.5
$ws.Range("J28").Value = 116.333336
$ws.Range("K28").Value = 242.5
$ws.Range("L28").Value = 116.333336
$ws.Range("M28").Value = 242.5
$ws.Range("N28").Value = -1086.333336
$ws.Range("H53").Value = 604.3043
$ws.Range("I53").Value = 245.38461
$ws.Range("K53").Value = 245.38461
$ws.Range("M53").Value = 391.61539
$ws.Range("H132").Value = 1676.2413
$ws.Range("I132").Value = 1116.92
$ws.Range("K132").Value = 3350.76
$ws.Range("M132").Value = -820.7600000000002
$ws.Range("H133").Value = 119590
$ws.Range("J133").Value = 119590
$ws.Range("L133").Value = 119590
$ws.Range("N133").Value = -129710

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 1750000
$ws.Range("J8").Value = 500000
$ws.Range("L8").Value = 500000
$ws.Range("N8").Value = -500288
$ws.Range("H36").Value = 1113473.9
$ws.Range("I36").Value = 2823.5715
$ws.Range("K36").Value = 2823.5715
$ws.Range("M36").Value = -2477.5715
$ws.Range("H61").Value = 10328.083
$ws.Range("I61").Value = 8358.666999999999
$ws.Range("K61").Value = 8358.666999999999
$ws.Range("M61").Value = -8146.666999999999
$ws.Range("H132").Value = 10976.412
$ws.Range("I132").Value = 7563.5454
$ws.Range("K132").Value = 22690.6362
$ws.Range("M132").Value = -20160.6362
$ws.Range("H136").Value = 10328.083
$ws.Range("I136").Value = 8358.666999999999
$ws.Range("K136").Value = 25076.001
$ws.Range("M136").Value = -22526.001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5204.488
$ws.Range("I134").Value = 4188.033
$ws.Range("K134").Value = 12564.099
$ws.Range("M134").Value = -10029.099

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13703572
$ws.Range("I31").Value = 34485060
$ws.Range("J31").Value = 6681.2046
$ws.Range("K31").Value = 34485060
$ws.Range("L31").Value = 6681.2046
$ws.Range("M31").Value = -34484765
$ws.Range("N31").Value = -7271.2046
$ws.Range("H34").Value = 13703572
$ws.Range("I34").Value = 34485060
$ws.Range("J34").Value = 6681.2046
$ws.Range("K34").Value = 34485060
$ws.Range("L34").Value = 6681.2046
$ws.Range("M34").Value = -34484858
$ws.Range("N34").Value = -7085.2046
$ws.Range("H58").Value = 4592.9287
$ws.Range("I58").Value = 2876.3333
$ws.Range("K58").Value = 2876.3333
$ws.Range("M58").Value = -2673.3333
$ws.Range("H94").Value = 76292.78999999999
$ws.Range("I94").Value = 149317
$ws.Range("K94").Value = 149317
$ws.Range("M94").Value = -148866
$ws.Range("H132").Value = 24586.191
$ws.Range("I132").Value = 3131.5
$ws.Range("J132").Value = 72859.25
$ws.Range("K132").Value = 9394.5
$ws.Range("L132").Value = 218577.75
$ws.Range("M132").Value = -6864.5
$ws.Range("N132").Value = -223637.75
$ws.Range("H134").Value = 5003.879
$ws.Range("I134").Value = 4352.967
$ws.Range("J134").Value = 11513
$ws.Range("K134").Value = 13058.901
$ws.Range("L134").Value = 34539
$ws.Range("M134").Value = -10523.901
$ws.Range("N134").Value = -39609
$ws.Range("H136").Value = 4592.9287
$ws.Range("I136").Value = 2876.3333
$ws.Range("K136").Value = 8628.999899999999
$ws.Range("M136").Value = -6078.999899999999
$ws.Range("H141").Value = 240483.78
$ws.Range("J141").Value = 260127.58
$ws.Range("L141").Value = 260127.58
$ws.Range("N141").Value = -270487.58

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 562
$ws.Range("I17").Value = 451.85715
$ws.Range("J17").Value = 754.75
$ws.Range("K17").Value = 1355.57145
$ws.Range("L17").Value = 2264.25
$ws.Range("M17").Value = -1186.57145
$ws.Range("N17").Value = -2602.25
$ws.Range("H113").Value = 2236.7058
$ws.Range("I113").Value = 1484.25
$ws.Range("J113").Value = 2905.5557
$ws.Range("K113").Value = 4452.75
$ws.Range("L113").Value = 8716.667099999999
$ws.Range("M113").Value = -2282.75
$ws.Range("N113").Value = -13056.6671
$ws.Range("H138").Value = 771781.4
$ws.Range("I138").Value = 1909.2
$ws.Range("J138").Value = 3338022
$ws.Range("K138").Value = 5727.6
$ws.Range("L138").Value = 10014066
$ws.Range("M138").Value = -587.6000000000004
$ws.Range("N138").Value = -10024346
$ws.Range("H139").Value = 2696.8333
$ws.Range("I139").Value = 1761.2858
$ws.Range("J139").Value = 4006.6
$ws.Range("K139").Value = 5283.857400000001
$ws.Range("L139").Value = 12019.8
$ws.Range("M139").Value = -143.8574000000008
$ws.Range("N139").Value = -22299.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 88494.5
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 88494.5
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H126").Value = 5110
$ws.Range("I126").Value = 7000
$ws.Range("J126").Value = 4873.75
$ws.Range("K126").Value = 21000
$ws.Range("L126").Value = 14621.25
$ws.Range("M126").Value = -18530
$ws.Range("N126").Value = -19561.25
$ws.Range("H132").Value = 6072.125
$ws.Range("I132").Value = 4391.1055
$ws.Range("J132").Value = 12460
$ws.Range("K132").Value = 13173.3165
$ws.Range("L132").Value = 37380
$ws.Range("M132").Value = -10643.3165
$ws.Range("N132").Value = -42440

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 544.38464
$ws.Range("I55").Value = 586.1111
$ws.Range("J55").Value = 450.5
$ws.Range("K55").Value = 586.1111
$ws.Range("L55").Value = 450.5
$ws.Range("M55").Value = -413.1111
$ws.Range("N55").Value = -796.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 35578.285
$ws.Range("I61").Value = 35000
$ws.Range("J61").Value = 35809.6
$ws.Range("K61").Value = 35000
$ws.Range("L61").Value = 35809.6
$ws.Range("M61").Value = -34708
$ws.Range("N61").Value = -36393.6
$ws.Range("H81").Value = 11905.596
$ws.Range("J81").Value = 14930.827
$ws.Range("L81").Value = 29861.654
$ws.Range("N81").Value = -31983.654
$ws.Range("H84").Value = 11905.596
$ws.Range("J84").Value = 14930.827
$ws.Range("L84").Value = 149308.27
$ws.Range("N84").Value = -159916.27
$ws.Range("H92").Value = 28916.166
$ws.Range("J92").Value = 28916.166
$ws.Range("L92").Value = 28916.166
$ws.Range("N92").Value = -33908.166
$ws.Range("H132").Value = 4738.6763
$ws.Range("I132").Value = 4480
$ws.Range("J132").Value = 5108.2144
$ws.Range("K132").Value = 13440
$ws.Range("L132").Value = 15324.6432
$ws.Range("M132").Value = -10910
$ws.Range("N132").Value = -20384.6432
